$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of old "Indicator" labels -> new labels (2021 PPP poverty-line rename)
$replacements = @{
    "Poverty $3.00 usd" = "Poverty $3.00/day (2021 PPP)"
    "Poverty $4.20 usd" = "Poverty $4.20/day (2021 PPP)"
    "Poverty $8.30 usd" = "Poverty $8.30/day (2021 PPP)"
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
